$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the spelling of "Purple halanchoe" -> "Purple kalanchoe" (name column for row 15)
$ws.Range("B15").Value = "Purple kalanchoe"

# Set explicit width for column B (name column)
$ws.Columns.Item(2).ColumnWidth = 48.5

# Update the selected/active cell to B15
$ws.Range("B15").Select()
